$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = [double]"0.0004098448753356933"
$ws.Range("M2").Value = [double]"0.08551836013793945"
$ws.Range("N2").Value = [double]"9.393692016601562e-05"
$ws.Range("L3").Value = [double]"8.559322357177735e-05"
$ws.Range("M3").Value = [double]"0.003118276596069336"
$ws.Range("N3").Value = [double]"5.197525024414062e-05"
$ws.Range("L4").Value = [double]"0.0003015847206115723"
$ws.Range("M4").Value = [double]"0.003179311752319336"
$ws.Range("N4").Value = [double]"0.0001177787780761719"
$ws.Range("L5").Value = [double]"7.24029541015625e-05"
$ws.Range("M5").Value = [double]"0.0007503032684326172"
$ws.Range("N5").Value = [double]"4.720687866210938e-05"
$ws.Range("E6").Value = [double]"0.9749791492910759"
$ws.Range("F6").Value = [double]"0.8958333333333334"
$ws.Range("G6").Value = [double]"0.8958333333333334"
$ws.Range("H6").Value = [double]"1040"
$ws.Range("I6").Value = [double]"15"
$ws.Range("J6").Value = [double]"15"
$ws.Range("K6").Value = [double]"129"
$ws.Range("L6").Value = [double]"0.008928798437118531"
$ws.Range("M6").Value = [double]"0.04216599464416504"
$ws.Range("N6").Value = [double]"0.007203340530395508"
$ws.Range("E7").Value = [double]"0.97581317764804"
$ws.Range("F7").Value = [double]"0.9020979020979021"
$ws.Range("H7").Value = [double]"1041"
$ws.Range("I7").Value = [double]"14"
$ws.Range("L7").Value = [double]"0.01020641756057739"
$ws.Range("M7").Value = [double]"0.1040611267089844"
$ws.Range("N7").Value = [double]"0.006453037261962891"
$ws.Range("L8").Value = [double]"0.0007186317443847656"
$ws.Range("M8").Value = [double]"0.009264707565307617"
$ws.Range("N8").Value = [double]"0.0002360343933105469"
$ws.Range("L9").Value = [double]"0.0001060323715209961"
$ws.Range("M9").Value = [double]"0.0005729198455810547"
$ws.Range("N9").Value = [double]"8.296966552734375e-05"
$ws.Range("L10").Value = [double]"0.0005566618442535401"
$ws.Range("M10").Value = [double]"0.003770351409912109"
$ws.Range("N10").Value = [double]"0.0001456737518310547"
$ws.Range("L11").Value = [double]"0.0001553776264190674"
$ws.Range("M11").Value = [double]"0.002080917358398438"
$ws.Range("N11").Value = [double]"9.1552734375e-05"
$ws.Range("L12").Value = [double]"0.002225195169448853"
$ws.Range("M12").Value = [double]"0.01387214660644531"
$ws.Range("N12").Value = [double]"0.0006763935089111328"
$ws.Range("L13").Value = [double]"0.0004849436283111572"
$ws.Range("M13").Value = [double]"0.001590251922607422"
$ws.Range("N13").Value = [double]"0.0003821849822998047"
$ws.Range("L14").Value = [double]"0.0002855598926544189"
$ws.Range("M14").Value = [double]"0.01058053970336914"
$ws.Range("N14").Value = [double]"9.846687316894531e-05"
$ws.Range("L15").Value = [double]"0.0001470100879669189"
$ws.Range("M15").Value = [double]"0.01251840591430664"
$ws.Range("N15").Value = [double]"4.553794860839844e-05"
$ws.Range("L16").Value = [double]"0.000534229040145874"
$ws.Range("M16").Value = [double]"0.007826328277587891"
$ws.Range("N16").Value = [double]"0.0002019405364990234"
$ws.Range("L17").Value = [double]"0.0001687929630279541"
$ws.Range("M17").Value = [double]"0.002122402191162109"
$ws.Range("N17").Value = [double]"0.0001037120819091797"
